$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date as literal text (matching the source data, which stores
# dates as plain strings rather than Excel date serials). Temporarily force
# a text number format so Excel's autoconvert doesn't turn the string into
# a date serial, then clear the formatting again so the new row matches the
# unformatted style of the existing data rows.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "09/05/2025"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = 603.0450000000001
$ws.Range("C4").Value = 0.08291255213126714
$ws.Range("D4").Value = 50
